$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5883970.5
$ws.Range("I33").Value = 16666752
$ws.Range("K33").Value = 16666752
$ws.Range("M33").Value = -16666523
$ws.Range("H64").Value = 9953.727999999999
$ws.Range("I64").Value = 3489.8
$ws.Range("J64").Value = 15340.333
$ws.Range("K64").Value = 3489.8
$ws.Range("L64").Value = 15340.333
$ws.Range("M64").Value = -3241.8
$ws.Range("N64").Value = -15836.333
$ws.Range("H67").Value = 9953.727999999999
$ws.Range("I67").Value = 3489.8
$ws.Range("J67").Value = 15340.333
$ws.Range("K67").Value = 3489.8
$ws.Range("L67").Value = 15340.333
$ws.Range("M67").Value = -2631.8
$ws.Range("N67").Value = -17056.333
$ws.Range("H87").Value = 59250
$ws.Range("J87").Value = 59250
$ws.Range("L87").Value = 59250
$ws.Range("N87").Value = -61746
$ws.Range("H90").Value = 59250
$ws.Range("J90").Value = 59250
$ws.Range("L90").Value = 177750
$ws.Range("N90").Value = -190230
$ws.Range("H98").Value = 1008.6667
$ws.Range("I98").Value = 509.75
$ws.Range("K98").Value = 509.75
$ws.Range("M98").Value = 988.25
$ws.Range("H107").Value = 4083.1333
$ws.Range("I107").Value = 2699.7778
$ws.Range("K107").Value = 2699.7778
$ws.Range("M107").Value = -779.7777999999998
$ws.Range("H112").Value = 605602.4
$ws.Range("J112").Value = 672758.5600000001
$ws.Range("L112").Value = 2018275.68
$ws.Range("N112").Value = -2020491.68
$ws.Range("H122").Value = 1008.6667
$ws.Range("I122").Value = 509.75
$ws.Range("K122").Value = 1529.25
$ws.Range("M122").Value = 920.75
$ws.Range("I137").Value = 56423.2
$ws.Range("J137").Value = 4366066
$ws.Range("K137").Value = 169269.6
$ws.Range("L137").Value = 13098198
$ws.Range("M137").Value = -166719.6
$ws.Range("N137").Value = -13103298

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3723.8647
$ws.Range("I32").Value = 3271.8484
$ws.Range("J32").Value = 7453
$ws.Range("K32").Value = 3271.8484
$ws.Range("L32").Value = 7453
$ws.Range("M32").Value = -2984.8484
$ws.Range("N32").Value = -8027
$ws.Range("H61").Value = 2710.389
$ws.Range("I61").Value = 2458.0588
$ws.Range("K61").Value = 2458.0588
$ws.Range("M61").Value = -2246.0588
$ws.Range("H136").Value = 2710.389
$ws.Range("I136").Value = 2458.0588
$ws.Range("K136").Value = 7374.176399999999
$ws.Range("M136").Value = -4824.176399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3029.543
$ws.Range("I86").Value = 1792.6786
$ws.Range("J86").Value = 7977
$ws.Range("K86").Value = 1792.6786
$ws.Range("L86").Value = 7977
$ws.Range("M86").Value = -669.6786
$ws.Range("N86").Value = -10223
$ws.Range("H89").Value = 3029.543
$ws.Range("I89").Value = 1792.6786
$ws.Range("J89").Value = 7977
$ws.Range("K89").Value = 8963.393
$ws.Range("L89").Value = 39885
$ws.Range("M89").Value = -3347.393
$ws.Range("N89").Value = -51117

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 25800
$ws.Range("J59").Value = 25800
$ws.Range("L59").Value = 25800
$ws.Range("N59").Value = -28090
$ws.Range("H86").Value = 7623.1
$ws.Range("I86").Value = 7622.875
$ws.Range("J86").Value = 7624
$ws.Range("K86").Value = 7622.875
$ws.Range("L86").Value = 7624
$ws.Range("M86").Value = -6499.875
$ws.Range("N86").Value = -9870
$ws.Range("H89").Value = 7623.1
$ws.Range("I89").Value = 7622.875
$ws.Range("J89").Value = 7624
$ws.Range("K89").Value = 38114.375
$ws.Range("L89").Value = 38120
$ws.Range("M89").Value = -32498.375
$ws.Range("N89").Value = -49352
$ws.Range("H99").Value = 3649.5
$ws.Range("I99").Value = 3699.3333
$ws.Range("K99").Value = 3699.3333
$ws.Range("M99").Value = -2201.3333
$ws.Range("H126").Value = 3649.5
$ws.Range("I126").Value = 3699.3333
$ws.Range("K126").Value = 11097.9999
$ws.Range("M126").Value = -8627.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 62.866665
$ws.Range("I38").Value = 51.666668
$ws.Range("K38").Value = 155.000004
$ws.Range("M38").Value = 191.999996
$ws.Range("H64").Value = 9998.75
$ws.Range("I64").Value = 9998.666999999999
$ws.Range("J64").Value = 9999
$ws.Range("K64").Value = 29996.001
$ws.Range("L64").Value = 29997
$ws.Range("M64").Value = -29726.001
$ws.Range("N64").Value = -30537
$ws.Range("H67").Value = 9998.75
$ws.Range("I67").Value = 9998.666999999999
$ws.Range("J67").Value = 9999
$ws.Range("K67").Value = 29996.001
$ws.Range("L67").Value = 29997
$ws.Range("M67").Value = -29060.001
$ws.Range("N67").Value = -31869
$ws.Range("H68").Value = 4168214.5
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 5001657.5
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 15004972.5
$ws.Range("N68").Value = -15006594.5
$ws.Range("M68").Value = -2189
$ws.Range("H71").Value = 4168214.5
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 5001657.5
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 45014917.5
$ws.Range("N71").Value = -45023029.5
$ws.Range("M71").Value = -4944

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12985.7
$ws.Range("I122").Value = 16729.572
$ws.Range("K122").Value = 50188.716
$ws.Range("M122").Value = -47738.716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 6644.75
$ws.Range("I42").Value = 6494.278
$ws.Range("K42").Value = 6494.278
$ws.Range("M42").Value = -5931.278
$ws.Range("H49").Value = 6644.75
$ws.Range("I49").Value = 6494.278
$ws.Range("K49").Value = 6494.278
$ws.Range("M49").Value = -6347.278
$ws.Range("H80").Value = 39849
$ws.Range("J80").Value = 39849
$ws.Range("L80").Value = 39849
$ws.Range("N80").Value = -42095
$ws.Range("H83").Value = 39849
$ws.Range("J83").Value = 39849
$ws.Range("L83").Value = 119547
$ws.Range("N83").Value = -130779

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H126").Value = 2877.75
$ws.Range("I126").Value = 2366.6365
$ws.Range("K126").Value = 7099.9095
$ws.Range("M126").Value = -4629.9095
